$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.259293079376221
$ws.Range("B1").Value = 2.510208368301392
$ws.Range("C1").Value = 4.730030059814453
$ws.Range("D1").Value = 2.01446533203125
$ws.Range("E1").Value = 1.151864886283875
